$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values for rows 2..96, regenerated from save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals"). Column G header is "K" (row 1).
$kValues = @(2,3,1,0,2,3,1,0,1,0,0,0,1,1,2,0,3,1,1,0,1,0,0,0,0,0,2,1,2,1,0,1,1,0,0,0,3,0,0,0,0,0,0,2,2,1,3,0,1,2,2,1,1,1,2,0,0,0,0,1,1,3,1,0,1,1,0,0,1,0,1,0,0,0,1,2,1,1,2,2,0,2,1,1,2,1,2,4,3,1,0,1,3,2,1)

$startRow = 2
for ($idx = 0; $idx -lt $kValues.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 7).Value = $kValues[$idx]
}
